$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A (shifts existing email/pass data to B/C)
$ws.Range("A1").EntireColumn.Insert()

# Insert a new header row at the top (shifts all data down by one row)
$ws.Range("A1").EntireRow.Insert()

# UUID values for each data row (2..11) -- fill these in BEFORE the header
# row so the shared-string table picks up the uuid values first (matching
# the original authoring order: data rows were entered, then headers).
$uuids = @(
    "9a757855-def8-4ffe-b2f7-882be9bada10",
    "9a735c1d-0131-4387-a82b-76ef177630d3",
    "9a735c1c-9c34-428a-bae3-acc6f1b094e8",
    "9a735c1c-4a50-4be1-a458-e22e360b165f",
    "9a735c1b-ebc4-4f7d-88cf-ca6796736b82",
    "9a735c1b-8e72-417f-875d-c3a59ddba5e2",
    "9a735c1b-37ee-4484-829d-f23803e1ccf0",
    "9a735c1a-d66a-4b34-b6c7-26983f935ea5",
    "9a735c1a-7714-42e1-a937-6540b962d3f0",
    "9a735c1a-182e-45b5-8041-3b9dee4e5b92"
)

for ($i = 0; $i -lt $uuids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $uuids[$i]
}

# Header row (written after the uuid values, per shared-string order)
$ws.Range("A1").Value = "uuid"
$ws.Range("B1").Value = "email"
$ws.Range("C1").Value = "pass"

# Column width for the new column A (columns B/C keep their original widths
# automatically, carried over by the column insert above).
$ws.Columns.Item(1).ColumnWidth = 40.3

# Update selection to match target state
$ws.Range("H8").Select()
